$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relational Model section ---
$ws.Range("B4").Value = 57
$ws.Range("B5").Value = 17
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = "2 enums done via JPA, which doesn't create a table"

# --- Object Model section ---
$ws.Range("B11").Value = 57
$ws.Range("B13").Value = 9
$ws.Range("B14").Value = 10
$ws.Range("B15").Value = 10
$ws.Range("B16").Value = 11
$ws.Range("B17").Value = 10
$ws.Range("B18").Value = 10

# --- User Interface library notes (entered ahead of other new comments,
#     matching the order the new shared strings were appended) ---
$ws.Range("C43").Value = "bootstrap theme and bootstrap"
$ws.Range("C42").Value = "(included jquery and jquery ui, but barely used - I started with javascript and retreated to jsp )as a plan B"

# --- Online Web Services note ---
$ws.Range("D29").Value = "(Note: Auth initial = fb login, but not ""real"" after that, so maybe that one doesn't count as much)"

# --- JWS or Express Web Services section ---
$ws.Range("C20").Value = "(most DAO CRUD methods are hooked up, but only one used in UI, so only counting that one)"
$ws.Range("C21").Value = "verifyAuthenticatedUser"
$ws.Range("C23").Value = "verifyAuthenticatedUser"

# --- Server Web Service Client section ---
$ws.Range("B30").Value = 3
$ws.Range("C31").Value = "Server side makes calls to Yelp and FlightAware in java, not Ajax.  UI does one ajax call to server ""verify"" method"
$ws.Range("B32").Value = 0
$ws.Range("B33").Value = 0
$ws.Range("B34").Value = 0

# --- User Interface section ---
$ws.Range("B37").Formula = "=22+6"
$ws.Range("B38").Formula = "=4+8"
$ws.Range("B39").Formula = "=1+1+4"
$ws.Range("B42").Value = 2
$ws.Range("B43").Value = 2

# --- View state: scroll / selection ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("C31").Select()
